# "Tried to implement Penality Reward System (unfinished)"
#
# Two targeted edits across the two sheets of the PO-analysis workbook:
#
# 1. "Weekly Quantity" sheet (sheet 1): remove the trailing week row
#    (A14 = 45130.99999999999 / B14 = 140), shrinking the used range
#    from A1:B14 down to A1:B13.
# 2. "Monthly Trend" sheet (sheet 2): the last month's requested
#    quantity (B7) drops from 620 to 480.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" -------------------------------------------
$wsWeekly = $wb.Worksheets.Item(1)

# Delete the entire 14th row (shifts nothing below it up since it's the
# last row) so the row disappears and the sheet's dimension becomes A1:B13.
$wsWeekly.Rows.Item(14).Delete()

# --- Sheet 2: "Monthly Trend" ----------------------------------------------
$wsMonthly = $wb.Worksheets.Item(2)

# Update the requested quantity for the last month (row 7, column B).
$wsMonthly.Cells.Item(7, 2).Value = 480
